$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete columns B (mainContact) and C (companyId), shifting D (projId) and E (refCompany) left.
$ws.Range("B1:C1").EntireColumn.Delete()

# Update the new B2 (formerly D2, projId row value) to the new content.
$ws.Range("B2").Value = "suporte Empresa 1"
